# Applies the cryptos.xlsx data refresh described by the diff:
# updated Price (D) / Volume(1h) (E) figures, plus the two row swaps
# (MinaProtocolToken <-> BinanceUSD at rows 30/31, and PaxDollar <-> FraxShare
# at rows 49/50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.034.41"
$ws.Range("E2").Value = "  +11.40%  "
$ws.Range("D3").Value = "1.810.73"
$ws.Range("E3").Value = "  +7.87%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.96"
$ws.Range("E5").Value = "  +3.83%  "
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.37"
$ws.Range("E8").Value = "  +4.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.81"
$ws.Range("E9").Value = "  +5.83%  "
$ws.Range("E10").Value = "  +6.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0666"
$ws.Range("E11").Value = "  +6.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0927"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "2.069.20"
$ws.Range("E13").Value = "  +7.70%  "
$ws.Range("D14").Value = "1.807.22"
$ws.Range("E14").Value = "  +8.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.638"
$ws.Range("E15").Value = "  +2.94%  "
$ws.Range("D16").Value = "33.950.07"
$ws.Range("E16").Value = "  +11.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.17"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.26"
$ws.Range("E18").Value = "  +7.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.28"
$ws.Range("E19").Value = "  +4.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "256.72"
$ws.Range("E20").Value = "  +5.01%  "
$ws.Range("E21").Value = "  +4.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.48"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.12"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.54"
$ws.Range("E27").Value = "  +4.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.10"
$ws.Range("E28").Value = "  +5.85%  "
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("B30").Value = "BinanceUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("B31").Value = "MinaProtocolToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.05"
$ws.Range("E31").Value = "  +398.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.86"
$ws.Range("E32").Value = "  +11.05%  "
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.20"
$ws.Range("E34").Value = "  +5.20%  "
$ws.Range("E35").Value = "  +6.45%  "
$ws.Range("D36").Value = "1.536.24"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("E38").Value = "  +4.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "84.21"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  +4.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.620"
$ws.Range("E41").Value = "  +5.30%  "
$ws.Range("E42").Value = "  +3.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  +8.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("E45").Value = "  +7.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0520"
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("D48").Value = "1.963.44"
$ws.Range("E48").Value = "  +7.96%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.71"
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.997"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.30"
$ws.Range("E51").Value = "  +2.15%  "
